$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Column D holds numeric-looking strings (e.g. "126.80", "1.40") that must
    # stay plain text (matching the inlineStr cells in the source file) instead
    # of being auto-coerced to numbers by Excel. Temporarily mark the cell as
    # Text, assign the value, then restore the original (Normal/General) style
    # so the cell formatting is left exactly as it was.
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '54.381.78'
$ws.Range('E2').Value = '  -2.89%  '
Set-TextValue 'D3' '2.282.38'
$ws.Range('E3').Value = '  -3.01%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.13%  '
Set-TextValue 'D5' '492.51'
$ws.Range('E5').Value = '  -2.41%  '
Set-TextValue 'D6' '126.80'
$ws.Range('E6').Value = '  -2.69%  '
$ws.Range('E8').Value = '  -1.82%  '
Set-TextValue 'D9' '2.281.53'
$ws.Range('E9').Value = '  -3.63%  '
Set-TextValue 'D10' '0.0942'
$ws.Range('E10').Value = '  -3.13%  '
Set-TextValue 'D11' '0.151'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('E12').Value = '  +0.19%  '
Set-TextValue 'D13' '4.62'
$ws.Range('E13').Value = '  -3.78%  '
Set-TextValue 'D14' '2.690.06'
$ws.Range('E14').Value = '  -2.90%  '
Set-TextValue 'D15' '21.49'
Set-TextValue 'D16' '54.294.74'
$ws.Range('E16').Value = '  -2.95%  '
Set-TextValue 'D18' '2.308.53'
$ws.Range('E18').Value = '  -2.14%  '
Set-TextValue 'D19' '9.97'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('E20').Value = '  +0.67%  '
Set-TextValue 'D21' '302.74'
$ws.Range('E21').Value = '  -2.71%  '
Set-TextValue 'D22' '6.44'
$ws.Range('E22').Value = '  +4.09%  '
$ws.Range('E23').Value = '  +0.28%  '
Set-TextValue 'D24' '5.37'
$ws.Range('E24').Value = '  -2.60%  '
Set-TextValue 'D25' '63.56'
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('E27').Value = '  +0.74%  '
Set-TextValue 'D28' '2.393.65'
$ws.Range('E28').Value = '  -2.90%  '
$ws.Range('E29').Value = '  +1.64%  '
Set-TextValue 'D30' '7.09'
$ws.Range('E30').Value = '  -0.76%  '
Set-TextValue 'D31' '168.54'
$ws.Range('E31').Value = '  -2.06%  '
Set-TextValue 'D32' '1.59'
$ws.Range('E32').Value = '  -2.81%  '
Set-TextValue 'D33' '0.0₃0683'
$ws.Range('E33').Value = '  -3.27%  '
$ws.Range('E34').Value = '  +2.08%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D36' '0.999'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D37' '1.08'
$ws.Range('E37').Value = '  +1.04%  '
Set-TextValue 'D38' '17.57'
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('E39').Value = '  +1.20%  '
Set-TextValue 'D40' '0.871'
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('E41').Value = '  -0.78%  '
Set-TextValue 'D42' '35.58'
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D43' '0.374'
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D44' '1.40'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -0.24%  '
Set-TextValue 'D46' '127.43'
$ws.Range('E46').Value = '  +1.68%  '
Set-TextValue 'D47' '4.84'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('E48').Value = '  -0.75%  '
Set-TextValue 'D49' '0.542'
$ws.Range('E49').Value = '  -2.84%  '
Set-TextValue 'D50' '238.80'
$ws.Range('E50').Value = '  -1.45%  '
Set-TextValue 'D51' '0.0478'
$ws.Range('E51').Value = '  -0.14%  '
